$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")
$ws.Activate()

# --- Day 5 block (rows 27-29) ---
$ws.Range("C27").Value = 7075
$ws.Range("C28").Value = 2700
$ws.Range("C29").Value = 2700

# --- Day 6 block (rows 33-35) ---
$ws.Range("C33").Value = 7075
$ws.Range("C34").Value = 2750
$ws.Range("C35").Value = 2750

# Re-touch the earlier day merges so they shuffle to the end of the
# mergeCells list (matches how Excel re-serialized them after this edit
# session).
$ws.Range("B2:C2").UnMerge()
$ws.Range("B2:C2").Merge()
$ws.Range("B8:C8").UnMerge()
$ws.Range("B8:C8").Merge()
$ws.Range("B14:C14").UnMerge()
$ws.Range("B14:C14").Merge()
$ws.Range("B20:C20").UnMerge()
$ws.Range("B20:C20").Merge()
$ws.Range("B26:C26").UnMerge()
$ws.Range("B26:C26").Merge()

# Update the view: scroll so row 18 is the top-left row shown, and move
# the active selection to the newly filled-in cell.
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C35").Select()
